# Commodity disaggregation rename: "Copper ores and concentrates" -> "Copper"
# The same shared text is used as the commodity label (column C) on every
# yearly worksheet (2000 .. 2100). Renaming one occurrence updates the
# shared-string table and every cell that referenced it.

$wb = $excel.ActiveWorkbook
$oldName = "Copper ores and concentrates"
$newName = "Copper"

foreach ($ws in $wb.Worksheets) {
    $found = $ws.Cells.Find($oldName)
    if ($found) {
        $firstAddress = $found.Address()
        do {
            $found.Value = $newName
            $found = $ws.Cells.FindNext($found)
        } while ($found -and ($found.Address() -ne $firstAddress))
    }
}
